$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = '2.623.77'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = "'593.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = "'166.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.52%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  -2.25%  '
$ws.Range("D9").Value = '2.623.90'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").Value = '66.908.77'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").Value = "'12.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.46%  '
$ws.Range("D20").Value = "'8.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.15%  '
$ws.Range("D21").Value = "'355.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("E23").Value = '  -2.92%  '
$ws.Range("E24").Value = '  +7.91%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -5.89%  '
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("D28").Value = '2.761.73'
$ws.Range("E28").Value = '  -1.98%  '
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("E30").Value = '  -1.54%  '
$ws.Range("D31").Value = "'547.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("E33").Value = '  -2.48%  '
$ws.Range("E34").Value = '  -1.40%  '
$ws.Range("E35").Value = '  +5.06%  '
$ws.Range("E37").Value = '  -4.91%  '
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = "'5.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.19%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'1.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.52%  '
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = "'40.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("E46").Value = '  -5.03%  '
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = "'151.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = "'0.577"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("E51").Value = '  -0.73%  '
